# Auto-generated script applying the scheduled-runner price/profit update
# to the Sheets/Zeromus_Profits.xlsx workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H15").Value = 2671.7576
$ws.Range("I15").Value = 2671.7576
$ws.Range("K15").Value = 8015.2728
$ws.Range("M15").Value = -7846.2728
$ws.Range("H76").Value = 4211.1113
$ws.Range("I76").Value = 3600
$ws.Range("J76").Value = 4700
$ws.Range("K76").Value = 3600
$ws.Range("L76").Value = 4700
$ws.Range("M76").Value = -3285
$ws.Range("N76").Value = -5330
$ws.Range("H79").Value = 4211.1113
$ws.Range("I79").Value = 3600
$ws.Range("J79").Value = 4700
$ws.Range("K79").Value = 3600
$ws.Range("L79").Value = 4700
$ws.Range("M79").Value = -2508
$ws.Range("N79").Value = -6884
$ws.Range("H106").Value = 4754.4443
$ws.Range("I106").Value = 4754.4443
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4754.4443
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents() | Out-Null
$ws.Range("N106").Value = -4123.4443
$ws.Range("H125").Value = 863.55554
$ws.Range("I125").Value = 795.3333
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 7157.9997
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -4697.9997
$ws.Range("N125").Value = -13920
$ws.Range("H132").Value = 7132.727
$ws.Range("I132").Value = 8660
$ws.Range("K132").Value = 25980
$ws.Range("M132").Value = -23450

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 11178.707
$ws.Range("I32").Value = 3286.1372
$ws.Range("J32").Value = 27950.416
$ws.Range("K32").Value = 3286.1372
$ws.Range("L32").Value = 27950.416
$ws.Range("M32").Value = -2999.1372
$ws.Range("N32").Value = -28524.416
$ws.Range("H61").Value = 1568.5385
$ws.Range("I61").Value = 1213.875
$ws.Range("K61").Value = 1213.875
$ws.Range("M61").Value = -1001.875
$ws.Range("H74").Value = 10976953
$ws.Range("I74").Value = 19565956
$ws.Range("J74").Value = 2115.8333
$ws.Range("K74").Value = 19565956
$ws.Range("L74").Value = 2115.8333
$ws.Range("M74").Value = -19565082
$ws.Range("N74").Value = -3863.8333
$ws.Range("H77").Value = 10976953
$ws.Range("I77").Value = 19565956
$ws.Range("J77").Value = 2115.8333
$ws.Range("K77").Value = 97829780
$ws.Range("L77").Value = 10579.1665
$ws.Range("M77").Value = -97825412
$ws.Range("N77").Value = -19315.1665
$ws.Range("H122").Value = 1492.5358
$ws.Range("I122").Value = 1220.4584
$ws.Range("K122").Value = 3661.3752
$ws.Range("M122").Value = -1211.3752
$ws.Range("H136").Value = 1568.5385
$ws.Range("I136").Value = 1213.875
$ws.Range("K136").Value = 3641.625
$ws.Range("M136").Value = -1091.625

$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 9151.305
$ws.Range("I94").Value = 333
$ws.Range("J94").Value = 29307.428
$ws.Range("K94").Value = 333
$ws.Range("L94").Value = 29307.428
$ws.Range("M94").Value = 118
$ws.Range("N94").Value = -30209.428
$ws.Range("H134").Value = 2718.862
$ws.Range("I134").Value = 2180.7
$ws.Range("J134").Value = 3914.7778
$ws.Range("K134").Value = 6542.099999999999
$ws.Range("L134").Value = 11744.3334
$ws.Range("M134").Value = -4007.099999999999
$ws.Range("N134").Value = -16814.3334

$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 2653.238
$ws.Range("I16").Value = 2578
$ws.Range("J16").Value = 2803.7144
$ws.Range("K16").Value = 2578
$ws.Range("L16").Value = 2803.7144
$ws.Range("M16").Value = -2291
$ws.Range("N16").Value = -3377.7144
$ws.Range("H31").Value = 5898136
$ws.Range("I31").Value = 5290375.5
$ws.Range("J31").Value = 6667965.5
$ws.Range("K31").Value = 5290375.5
$ws.Range("L31").Value = 6667965.5
$ws.Range("M31").Value = -5290080.5
$ws.Range("N31").Value = -6668555.5
$ws.Range("H34").Value = 5898136
$ws.Range("I34").Value = 5290375.5
$ws.Range("J34").Value = 6667965.5
$ws.Range("K34").Value = 5290375.5
$ws.Range("L34").Value = 6667965.5
$ws.Range("M34").Value = -5290173.5
$ws.Range("N34").Value = -6668369.5
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 11000
$ws.Range("J54").Value = 19000
$ws.Range("K54").Value = 11000
$ws.Range("L54").Value = 19000
$ws.Range("M54").Value = -10342
$ws.Range("N54").Value = -20316
$ws.Range("H58").Value = 1818.875
$ws.Range("I58").Value = 1179.5714
$ws.Range("K58").Value = 1179.5714
$ws.Range("M58").Value = -976.5714
$ws.Range("H99").Value = 13111.25
$ws.Range("I99").Value = 13148.333
$ws.Range("J99").Value = 13000
$ws.Range("K99").Value = 13148.333
$ws.Range("L99").Value = 13000
$ws.Range("M99").Value = -11650.333
$ws.Range("N99").Value = -15996
$ws.Range("H100").Value = 36466.668
$ws.Range("J100").Value = 36466.668
$ws.Range("L100").Value = 36466.668
$ws.Range("N100").Value = -38630.668
$ws.Range("H113").Value = 2653.238
$ws.Range("I113").Value = 2578
$ws.Range("J113").Value = 2803.7144
$ws.Range("K113").Value = 2578
$ws.Range("L113").Value = 2803.7144
$ws.Range("M113").Value = -408
$ws.Range("N113").Value = -7143.7144
$ws.Range("H122").Value = 995
$ws.Range("I122").Value = 862.1429000000001
$ws.Range("J122").Value = 1460
$ws.Range("K122").Value = 2586.4287
$ws.Range("L122").Value = 4380
$ws.Range("M122").Value = -136.4287000000004
$ws.Range("N122").Value = -9280
$ws.Range("H126").Value = 13111.25
$ws.Range("I126").Value = 13148.333
$ws.Range("J126").Value = 13000
$ws.Range("K126").Value = 39444.999
$ws.Range("L126").Value = 39000
$ws.Range("M126").Value = -36974.999
$ws.Range("N126").Value = -43940
$ws.Range("H132").Value = 2049.625
$ws.Range("I132").Value = 1397.4286
$ws.Range("J132").Value = 2556.889
$ws.Range("K132").Value = 4192.2858
$ws.Range("L132").Value = 7670.667
$ws.Range("M132").Value = -1662.2858
$ws.Range("N132").Value = -12730.667
$ws.Range("H134").Value = 4300.095
$ws.Range("I134").Value = 4605.933
$ws.Range("J134").Value = 3535.5
$ws.Range("K134").Value = 13817.799
$ws.Range("L134").Value = 10606.5
$ws.Range("M134").Value = -11282.799
$ws.Range("N134").Value = -15676.5
$ws.Range("H136").Value = 1818.875
$ws.Range("I136").Value = 1179.5714
$ws.Range("K136").Value = 3538.7142
$ws.Range("M136").Value = -988.7142000000003

$ws = $wb.Worksheets.Item(5)
$ws.Range("H37").Value = 62730.453
$ws.Range("J37").Value = 62730.453
$ws.Range("L37").Value = 188191.359
$ws.Range("N37").Value = -188415.359
$ws.Range("H137").Value = 3976.3057
$ws.Range("J137").Value = 3973.2424
$ws.Range("L137").Value = 11919.7272
$ws.Range("N137").Value = -22119.7272
$ws.Range("H140").Value = 1077.75
$ws.Range("I140").Value = 915
$ws.Range("J140").Value = 2000
$ws.Range("K140").Value = 2745
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = 2435
$ws.Range("N140").Value = -16360

$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 37.8125
$ws.Range("I2").Value = 39.833332
$ws.Range("J2").Value = 36.6
$ws.Range("K2").Value = 39.833332
$ws.Range("L2").Value = 36.6
$ws.Range("M2").Value = 73.166668
$ws.Range("N2").Value = -262.6
$ws.Range("H80").Value = 3058.1538
$ws.Range("J80").Value = 2955.6
$ws.Range("L80").Value = 2955.6
$ws.Range("N80").Value = -4951.6
$ws.Range("H83").Value = 3058.1538
$ws.Range("J83").Value = 2955.6
$ws.Range("L83").Value = 14778
$ws.Range("N83").Value = -24762
$ws.Range("H97").Value = 2303.3333
$ws.Range("I97").Value = 2303.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2303.3333
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents() | Out-Null
$ws.Range("N97").Value = -1807.3333
$ws.Range("H126").Value = 2510.5264
$ws.Range("I126").Value = 2763.6365
$ws.Range("J126").Value = 2162.5
$ws.Range("K126").Value = 8290.9095
$ws.Range("L126").Value = 6487.5
$ws.Range("M126").Value = -5820.9095
$ws.Range("N126").Value = -11427.5
$ws.Range("H131").Value = 56700
$ws.Range("J131").Value = 56700
$ws.Range("L131").Value = 56700
$ws.Range("N131").Value = -66780
$ws.Range("H132").Value = 2937.2222
$ws.Range("I132").Value = 2190.7273
$ws.Range("J132").Value = 4110.2856
$ws.Range("K132").Value = 6572.1819
$ws.Range("L132").Value = 12330.8568
$ws.Range("M132").Value = -4042.1819
$ws.Range("N132").Value = -17390.8568

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 2376
$ws.Range("I7").Value = 2002
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 2002
$ws.Range("L7").Value = 2750
$ws.Range("M7").Value = -1890
$ws.Range("N7").Value = -2974
$ws.Range("H126").Value = 2376
$ws.Range("I126").Value = 2002
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 6006
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -3536
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 22511132
$ws.Range("I132").Value = 56272896
$ws.Range("J132").Value = 3289.4666
$ws.Range("K132").Value = 168818688
$ws.Range("L132").Value = 9868.399800000001
$ws.Range("M132").Value = -168816158
$ws.Range("N132").Value = -14928.3998

$ws = $wb.Worksheets.Item(8)
$ws.Range("H64").Value = 18707.072
$ws.Range("J64").Value = 18707.072
$ws.Range("L64").Value = 18707.072
$ws.Range("N64").Value = -19203.072
$ws.Range("H67").Value = 18707.072
$ws.Range("J67").Value = 18707.072
$ws.Range("L67").Value = 18707.072
$ws.Range("N67").Value = -20423.072
$ws.Range("H81").Value = 50002544
$ws.Range("I81").Value = 76925100
$ws.Range("J81").Value = 3500.8572
$ws.Range("K81").Value = 153850200
$ws.Range("L81").Value = 7001.7144
$ws.Range("M81").Value = -153849139
$ws.Range("N81").Value = -9123.714400000001
$ws.Range("H84").Value = 50002544
$ws.Range("I84").Value = 76925100
$ws.Range("J84").Value = 3500.8572
$ws.Range("K84").Value = 769251000
$ws.Range("L84").Value = 35008.572
$ws.Range("M84").Value = -769245696
$ws.Range("N84").Value = -45616.572
$ws.Range("H126").Value = 7168
$ws.Range("I126").Value = 8401.6
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 25204.8
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -22734.8
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 1590.9166
$ws.Range("I132").Value = 831.0714
$ws.Range("J132").Value = 4250.375
$ws.Range("K132").Value = 2493.2142
$ws.Range("L132").Value = 12751.125
$ws.Range("M132").Value = 36.78579999999965
$ws.Range("N132").Value = -17811.125
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 781.8182
$ws.Range("J136").Value = 2816.6667
$ws.Range("K136").Value = 2345.4546
$ws.Range("L136").Value = 8450.000100000001
$ws.Range("M136").Value = 204.5454
$ws.Range("N136").Value = -13550.0001
